# Deploying to gh-pages from @ LinuxForHealth/alvearie-fhir-ig@80fa500adfae01c9a5dd7ef65e90accc96781b5c 🚀
#
# Updates the StructureDefinition-employee-gender workbook:
#  - Metadata sheet: URL / Version / Date / Publisher moved from the
#    "Alvearie" / ibm.com identity to the "LinuxForHealth" identity.
#  - Elements sheet: the fixed value of Extension.url (same URL string)
#    is updated to match, and the root "Extension" row's rolled-up
#    Constraint(s) text is cleared (it now lives solely on the
#    Extension.extension row).

$wb = $excel.ActiveWorkbook

$newUrl = "http://linuxforhealth.org/fhir/cdm/StructureDefinition/employee-gender"

$meta = $wb.Worksheets.Item("Metadata")
$meta.Range("B2").Value = $newUrl
$meta.Range("B3").Value = "8.0.0"
$meta.Range("B8").Value = "2022-11-10T16:00:46+00:00"
$meta.Range("B9").Value = "LinuxForHealth Team"

$elements = $wb.Worksheets.Item("Elements")
# Extension.url (row 5) Fixed Value column (Q) carries the same URL text.
$elements.Range("Q5").Value = $newUrl
# Root Extension row (row 2) Constraint(s) column (AI) no longer carries
# the rolled-up ele-1/ext-1 constraint text.
$elements.Range("AI2").Value = ""
